$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.759382333333334
$ws.Range("H2").Value = 23.278147
$ws.Range("I2").Value = 0.03884312320086166
$ws.Range("J2").Value = 0.04014625174390325
$ws.Range("M2").Value = 2.035934
$ws.Range("N2").Value = 6.107802
$ws.Range("O2").Value = 0.03126880699186227
$ws.Range("P2").Value = 0.04430738339814538
$ws.Range("Q2").Value = 15.79759031143267
$ws.Range("R2").Value = 142.178312802894
$ws.Range("S2").Value = 0.001214578122328871
$ws.Range("T2").Value = 0.001778775368015584
$ws.Range("G3").Value = 7.759382333333334
$ws.Range("H3").Value = 23.278147
$ws.Range("I3").Value = 0.03884312320086166
$ws.Range("J3").Value = 0.04014625174390325
$ws.Range("O3").Value = 0.0582841555718936
$ws.Range("P3").Value = 0.08258768643246023
$ws.Range("Q3").Value = 29.44625330963878
$ws.Range("R3").Value = 265.016279786749
$ws.Range("S3").Value = 0.00226393863553725
$ws.Range("T3").Value = 0.003315586050464091
$ws.Range("G4").Value = 7.759382333333334
$ws.Range("H4").Value = 23.278147
$ws.Range("I4").Value = 0.03884312320086166
$ws.Range("J4").Value = 0.04014625174390325
$ws.Range("M4").Value = 0.451002
$ws.Range("N4").Value = 1.353006
$ws.Range("O4").Value = 0.006926695310822388
$ws.Range("P4").Value = 0.009815012926416261
$ws.Range("Q4").Value = 3.499496951098
$ws.Range("R4").Value = 31.495472559882
$ws.Range("S4").Value = 0.0002690544793331047
$ws.Range("T4").Value = 0.0003940359798135718
$ws.Range("G5").Value = 7.759382333333334
$ws.Range("H5").Value = 23.278147
$ws.Range("I5").Value = 0.03884312320086166
$ws.Range("J5").Value = 0.04014625174390325
$ws.Range("M5").Value = 57.4814495
$ws.Range("N5").Value = 114.962899
$ws.Range("O5").Value = 0.8828264325012393
$ws.Range("P5").Value = 0.8339669888701803
$ws.Range("Q5").Value = 446.0205437446922
$ws.Range("R5").Value = 2676.123262468153
$ws.Range("S5").Value = 0.03429173588262282
$ws.Range("T5").Value = 0.03348064868128722
$ws.Range("G6").Value = 7.759382333333334
$ws.Range("H6").Value = 23.278147
$ws.Range("I6").Value = 0.03884312320086166
$ws.Range("J6").Value = 0.04014625174390325
$ws.Range("M6").Value = 1.347395
$ws.Range("N6").Value = 4.042185
$ws.Range("O6").Value = 0.02069390962418245
$ws.Range("P6").Value = 0.02932292837279799
$ws.Range("Q6").Value = 10.45495295902167
$ws.Range("R6").Value = 94.094576631195
$ws.Range("S6").Value = 0.0008038160810396155
$ws.Range("T6").Value = 0.001177205664322791
$ws.Range("I7").Value = 0.8631909770948131
$ws.Range("J7").Value = 0.8921497401307179
$ws.Range("M7").Value = 2.035934
$ws.Range("N7").Value = 6.107802
$ws.Range("O7").Value = 0.03126880699186227
$ws.Range("P7").Value = 0.04430738339814538
$ws.Range("Q7").Value = 351.0618172013166
$ws.Range("R7").Value = 3159.556354811849
$ws.Range("S7").Value = 0.02699095205989472
$ws.Range("T7").Value = 0.03952882058452749
$ws.Range("I8").Value = 0.8631909770948131
$ws.Range("J8").Value = 0.8921497401307179
$ws.Range("O8").Value = 0.0582841555718936
$ws.Range("P8").Value = 0.08258768643246023
$ws.Range("S8").Value = 0.05031035719724893
$ws.Range("T8").Value = 0.0736805829887166
$ws.Range("I9").Value = 0.8631909770948131
$ws.Range("J9").Value = 0.8921497401307179
$ws.Range("M9").Value = 0.451002
$ws.Range("N9").Value = 1.353006
$ws.Range("O9").Value = 0.006926695310822388
$ws.Range("P9").Value = 0.009815012926416261
$ws.Range("Q9").Value = 77.76754142394999
$ws.Range("R9").Value = 699.9078728155499
$ws.Range("S9").Value = 0.005979060893386837
$ws.Range("T9").Value = 0.008756461231681905
$ws.Range("I10").Value = 0.8631909770948131
$ws.Range("J10").Value = 0.8921497401307179
$ws.Range("M10").Value = 57.4814495
$ws.Range("N10").Value = 114.962899
$ws.Range("O10").Value = 0.8828264325012393
$ws.Range("P10").Value = 0.8339669888701803
$ws.Range("Q10").Value = 9911.687764355678
$ws.Range("R10").Value = 59470.12658613407
$ws.Range("S10").Value = 0.7620478108758728
$ws.Range("T10").Value = 0.7440234323981286
$ws.Range("I11").Value = 0.8631909770948131
$ws.Range("J11").Value = 0.8921497401307179
$ws.Range("M11").Value = 1.347395
$ws.Range("N11").Value = 4.042185
$ws.Range("O11").Value = 0.02069390962418245
$ws.Range("P11").Value = 0.02932292837279799
$ws.Range("Q11").Value = 232.3351037842916
$ws.Range("R11").Value = 2091.015934058625
$ws.Range("S11").Value = 0.0178627960684098
$ws.Range("T11").Value = 0.02616044292766338
$ws.Range("G12").Value = 0.05240566666666666
$ws.Range("H12").Value = 0.157217
$ws.Range("I12").Value = 0.0002623404388789996
$ws.Range("J12").Value = 0.0002711415672571033
$ws.Range("M12").Value = 2.035934
$ws.Range("N12").Value = 6.107802
$ws.Range("O12").Value = 0.03126880699186227
$ws.Range("P12").Value = 0.04430738339814538
$ws.Range("Q12").Value = 0.1066944785593333
$ws.Range("R12").Value = 0.9602503070339999
$ws.Range("S12").Value = 0.000008203072549467877
$ws.Range("T12").Value = 0.0000120135733756345
$ws.Range("G13").Value = 0.05240566666666666
$ws.Range("H13").Value = 0.157217
$ws.Range("I13").Value = 0.0002623404388789996
$ws.Range("J13").Value = 0.0002711415672571033
$ws.Range("O13").Value = 0.0582841555718936
$ws.Range("P13").Value = 0.08258768643246023
$ws.Range("Q13").Value = 0.1988754348265556
$ws.Range("R13").Value = 1.789878913439
$ws.Range("S13").Value = 0.00001529029095242245
$ws.Range("T13").Value = 0.00002239295473543547
$ws.Range("G14").Value = 0.05240566666666666
$ws.Range("H14").Value = 0.157217
$ws.Range("I14").Value = 0.0002623404388789996
$ws.Range("J14").Value = 0.0002711415672571033
$ws.Range("M14").Value = 0.451002
$ws.Range("N14").Value = 1.353006
$ws.Range("O14").Value = 0.006926695310822388
$ws.Range("P14").Value = 0.009815012926416261
$ws.Range("Q14").Value = 0.023635060478
$ws.Range("R14").Value = 0.212715544302
$ws.Range("S14").Value = 0.000001817152287822254
$ws.Range("T14").Value = 0.000002661257987517233
$ws.Range("G15").Value = 0.05240566666666666
$ws.Range("H15").Value = 0.157217
$ws.Range("I15").Value = 0.0002623404388789996
$ws.Range("J15").Value = 0.0002711415672571033
$ws.Range("M15").Value = 57.4814495
$ws.Range("N15").Value = 114.962899
$ws.Range("O15").Value = 0.8828264325012393
$ws.Range("P15").Value = 0.8339669888701803
$ws.Range("Q15").Value = 3.012353682013833
$ws.Range("R15").Value = 18.074122092083
$ws.Range("S15").Value = 0.0002316010737563566
$ws.Range("T15").Value = 0.0002261231164029479
$ws.Range("G16").Value = 0.05240566666666666
$ws.Range("H16").Value = 0.157217
$ws.Range("I16").Value = 0.0002623404388789996
$ws.Range("J16").Value = 0.0002711415672571033
$ws.Range("M16").Value = 1.347395
$ws.Range("N16").Value = 4.042185
$ws.Range("O16").Value = 0.02069390962418245
$ws.Range("P16").Value = 0.02932292837279799
$ws.Range("Q16").Value = 0.07061113323833332
$ws.Range("R16").Value = 0.635500199145
$ws.Range("S16").Value = 0.000005428849332930375
$ws.Range("T16").Value = 0.000007950664755568229
$ws.Range("G17").Value = 19.452549
$ws.Range("H17").Value = 38.905098
$ws.Range("I17").Value = 0.09737859599105524
$ws.Range("J17").Value = 0.06709700125311635
$ws.Range("M17").Value = 2.035934
$ws.Range("N17").Value = 6.107802
$ws.Range("O17").Value = 0.03126880699186227
$ws.Range("P17").Value = 0.04430738339814538
$ws.Range("Q17").Value = 39.60410589576599
$ws.Range("R17").Value = 237.624635374596
$ws.Range("S17").Value = 0.003044912523182839
$ws.Range("T17").Value = 0.002972892559387667
$ws.Range("G18").Value = 19.452549
$ws.Range("H18").Value = 38.905098
$ws.Range("I18").Value = 0.09737859599105524
$ws.Range("J18").Value = 0.06709700125311635
$ws.Range("O18").Value = 0.0582841555718936
$ws.Range("P18").Value = 0.08258768643246023
$ws.Range("Q18").Value = 73.820912640361
$ws.Range("R18").Value = 442.925475842166
$ws.Range("S18").Value = 0.005675629238115237
$ws.Range("T18").Value = 0.005541386100050764
$ws.Range("G19").Value = 19.452549
$ws.Range("H19").Value = 38.905098
$ws.Range("I19").Value = 0.09737859599105524
$ws.Range("J19").Value = 0.06709700125311635
$ws.Range("M19").Value = 0.451002
$ws.Range("N19").Value = 1.353006
$ws.Range("O19").Value = 0.006926695310822388
$ws.Range("P19").Value = 0.009815012926416261
$ws.Range("Q19").Value = 8.773138504097998
$ws.Range("R19").Value = 52.63883102458799
$ws.Range("S19").Value = 0.0006745118642257101
$ws.Range("T19").Value = 0.0006585579346231051
$ws.Range("G20").Value = 19.452549
$ws.Range("H20").Value = 38.905098
$ws.Range("I20").Value = 0.09737859599105524
$ws.Range("J20").Value = 0.06709700125311635
$ws.Range("M20").Value = 57.4814495
$ws.Range("N20").Value = 114.962899
$ws.Range("O20").Value = 0.8828264325012393
$ws.Range("P20").Value = 0.8339669888701803
$ws.Range("Q20").Value = 1118.160712989775
$ws.Range("R20").Value = 4472.642851959101
$ws.Range("S20").Value = 0.08596839850076278
$ws.Range("T20").Value = 0.05595668409728016
$ws.Range("G21").Value = 19.452549
$ws.Range("H21").Value = 38.905098
$ws.Range("I21").Value = 0.09737859599105524
$ws.Range("J21").Value = 0.06709700125311635
$ws.Range("M21").Value = 1.347395
$ws.Range("N21").Value = 4.042185
$ws.Range("O21").Value = 0.02069390962418245
$ws.Range("P21").Value = 0.02932292837279799
$ws.Range("Q21").Value = 26.21026725985499
$ws.Range("R21").Value = 157.26160355913
$ws.Range("S21").Value = 0.002015143864768672
$ws.Range("T21").Value = 0.001967480561774668
$ws.Range("G22").Value = 0.06491533333333334
$ws.Range("H22").Value = 0.194746
$ws.Range("I22").Value = 0.0003249632743909987
$ws.Range("J22").Value = 0.0003358653050055137
$ws.Range("M22").Value = 2.035934
$ws.Range("N22").Value = 6.107802
$ws.Range("O22").Value = 0.03126880699186227
$ws.Range("P22").Value = 0.04430738339814538
$ws.Range("Q22").Value = 0.1321633342546667
$ws.Range("R22").Value = 1.189470008292
$ws.Range("S22").Value = 0.00001016121390637572
$ws.Range("T22").Value = 0.00001488131283901433
$ws.Range("G23").Value = 0.06491533333333334
$ws.Range("H23").Value = 0.194746
$ws.Range("I23").Value = 0.0003249632743909987
$ws.Range("J23").Value = 0.0003358653050055137
$ws.Range("O23").Value = 0.0582841555718936
$ws.Range("P23").Value = 0.08258768643246023
$ws.Range("Q23").Value = 0.2463486482424445
$ws.Range("R23").Value = 2.217137834182
$ws.Range("S23").Value = 0.00001894021003975692
$ws.Range("T23").Value = 0.00002773833849333798
$ws.Range("G24").Value = 0.06491533333333334
$ws.Range("H24").Value = 0.194746
$ws.Range("I24").Value = 0.0003249632743909987
$ws.Range("J24").Value = 0.0003358653050055137
$ws.Range("M24").Value = 0.451002
$ws.Range("N24").Value = 1.353006
$ws.Range("O24").Value = 0.006926695310822388
$ws.Range("P24").Value = 0.009815012926416261
$ws.Range("Q24").Value = 0.029276945164
$ws.Range("R24").Value = 0.263492506476
$ws.Range("S24").Value = 0.00000225092158891362
$ws.Range("T24").Value = 0.000003296522310163857
$ws.Range("G25").Value = 0.06491533333333334
$ws.Range("H25").Value = 0.194746
$ws.Range("I25").Value = 0.0003249632743909987
$ws.Range("J25").Value = 0.0003358653050055137
$ws.Range("M25").Value = 57.4814495
$ws.Range("N25").Value = 114.962899
$ws.Range("O25").Value = 0.8828264325012393
$ws.Range("P25").Value = 0.8339669888701803
$ws.Range("Q25").Value = 3.731427454775667
$ws.Range("R25").Value = 22.388564728654
$ws.Range("S25").Value = 0.0002868861682245267
$ws.Range("T25").Value = 0.0002801005770814129
$ws.Range("G26").Value = 0.06491533333333334
$ws.Range("H26").Value = 0.194746
$ws.Range("I26").Value = 0.0003249632743909987
$ws.Range("J26").Value = 0.0003358653050055137
$ws.Range("M26").Value = 1.347395
$ws.Range("N26").Value = 4.042185
$ws.Range("O26").Value = 0.02069390962418245
$ws.Range("P26").Value = 0.02932292837279799
$ws.Range("Q26").Value = 0.08746659555666667
$ws.Range("R26").Value = 0.78719936001
$ws.Range("S26").Value = 0.000006724760631425729
$ws.Range("T26").Value = 0.000009848554281584626
